$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.878.68"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "'2.603.27"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'554.73"
$ws.Range("E5").Value = "  +3.27%  "
$ws.Range("D6").Value = "'143.60"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +4.20%  "
$ws.Range("D9").Value = "'6.83"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("E11").Value = "  +4.91%  "
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "'3.063.95"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "'58.837.12"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").Value = "'20.88"
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("D16").Value = "'2.632.33"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("E17").Value = "  -2.02%  "
$ws.Range("D18").Value = "'4.46"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "'337.30"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").Value = "'10.10"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").Value = "'6.16"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'66.55"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -2.47%  "
$ws.Range("D27").Value = "'7.19"
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'1.67"
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").Value = "'5.99"
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("D32").Value = "'154.08"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").Value = "'19.00"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").Value = "'3.94"
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("D35").Value = "'0.908"
$ws.Range("E35").Value = "  +8.19%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "'0.877"
$ws.Range("E36").Value = "  +4.64%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.13"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("D38").Value = "'37.05"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "'281.95"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").Value = "'0.997"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Value = "'0.0539"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "'0.0954"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'4.66"
$ws.Range("E48").Value = "  +1.88%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'1.947.64"
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("D50").Value = "'118.73"
$ws.Range("E50").Value = "  +6.52%  "
$ws.Range("D51").Value = "'17.94"
$ws.Range("E51").Value = "  -2.79%  "
